$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = 44669
$ws.Cells.Item(2,9).Value = "Primera"
$ws.Cells.Item(2,10).Value = 130
$ws.Cells.Item(2,11).Value = 4500
$ws.Cells.Item(2,12).Value = 5000
$ws.Cells.Item(2,13).Value = 4750
$ws.Cells.Item(2,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(2,16).Value = 79
$ws.Cells.Item(2,17).Value = 60

# Row 3
$ws.Cells.Item(3,4).Value = 44627
$ws.Cells.Item(3,9).Value = "Primera"
$ws.Cells.Item(3,10).Value = 120
$ws.Cells.Item(3,11).Value = 4000
$ws.Cells.Item(3,12).Value = 4500
$ws.Cells.Item(3,13).Value = 4250
$ws.Cells.Item(3,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(3,16).Value = 71
$ws.Cells.Item(3,17).Value = 60

# Row 4
$ws.Cells.Item(4,4).Value = 44785
$ws.Cells.Item(4,9).Value = "Primera"
$ws.Cells.Item(4,10).Value = 130
$ws.Cells.Item(4,11).Value = 7000
$ws.Cells.Item(4,12).Value = 8000
$ws.Cells.Item(4,13).Value = 7500
$ws.Cells.Item(4,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(4,16).Value = 125
$ws.Cells.Item(4,17).Value = 60

# Row 5
$ws.Cells.Item(5,4).Value = 44764
$ws.Cells.Item(5,9).Value = "Primera"
$ws.Cells.Item(5,10).Value = 120
$ws.Cells.Item(5,11).Value = 7000
$ws.Cells.Item(5,12).Value = 8000
$ws.Cells.Item(5,13).Value = 7500
$ws.Cells.Item(5,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(5,16).Value = 125
$ws.Cells.Item(5,17).Value = 60

# Row 6
$ws.Cells.Item(6,4).Value = 44400
$ws.Cells.Item(6,9).Value = "Primera"
$ws.Cells.Item(6,10).Value = 120
$ws.Cells.Item(6,11).Value = 9000
$ws.Cells.Item(6,12).Value = 10000
$ws.Cells.Item(6,13).Value = 9500
$ws.Cells.Item(6,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(6,16).Value = 158
$ws.Cells.Item(6,17).Value = 60

# Row 7
$ws.Cells.Item(7,4).Value = 44382
$ws.Cells.Item(7,9).Value = "Primera"
$ws.Cells.Item(7,10).Value = 160
$ws.Cells.Item(7,11).Value = 7000
$ws.Cells.Item(7,12).Value = 8000
$ws.Cells.Item(7,13).Value = 7438
$ws.Cells.Item(7,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7,16).Value = 124
$ws.Cells.Item(7,17).Value = 60

# Row 8
$ws.Cells.Item(8,4).Value = 44760
$ws.Cells.Item(8,9).Value = "Primera"
$ws.Cells.Item(8,10).Value = 130
$ws.Cells.Item(8,11).Value = 7000
$ws.Cells.Item(8,12).Value = 7500
$ws.Cells.Item(8,13).Value = 7250
$ws.Cells.Item(8,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(8,16).Value = 121
$ws.Cells.Item(8,17).Value = 60

# Row 9
$ws.Cells.Item(9,4).Value = 44362
$ws.Cells.Item(9,9).Value = "Primera"
$ws.Cells.Item(9,10).Value = 120
$ws.Cells.Item(9,11).Value = 8000
$ws.Cells.Item(9,12).Value = 9000
$ws.Cells.Item(9,13).Value = 8500
$ws.Cells.Item(9,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(9,16).Value = 142
$ws.Cells.Item(9,17).Value = 60

# Row 10
$ws.Cells.Item(10,4).Value = 44963
$ws.Cells.Item(10,9).Value = "Primera"
$ws.Cells.Item(10,10).Value = 130
$ws.Cells.Item(10,11).Value = 4000
$ws.Cells.Item(10,12).Value = 4500
$ws.Cells.Item(10,13).Value = 4250
$ws.Cells.Item(10,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(10,16).Value = 71
$ws.Cells.Item(10,17).Value = 60

# Row 11
$ws.Cells.Item(11,4).Value = 44740
$ws.Cells.Item(11,9).Value = "Primera"
$ws.Cells.Item(11,10).Value = 120
$ws.Cells.Item(11,11).Value = 6000
$ws.Cells.Item(11,12).Value = 7000
$ws.Cells.Item(11,13).Value = 6500
$ws.Cells.Item(11,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(11,16).Value = 108
$ws.Cells.Item(11,17).Value = 60

# Row 12
$ws.Cells.Item(12,4).Value = 44281
$ws.Cells.Item(12,9).Value = "Primera"
$ws.Cells.Item(12,10).Value = 120
$ws.Cells.Item(12,11).Value = 5500
$ws.Cells.Item(12,12).Value = 6000
$ws.Cells.Item(12,13).Value = 5750
$ws.Cells.Item(12,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(12,16).Value = 96
$ws.Cells.Item(12,17).Value = 60

# Row 13
$ws.Cells.Item(13,4).Value = 44827
$ws.Cells.Item(13,9).Value = "Primera"
$ws.Cells.Item(13,10).Value = 120
$ws.Cells.Item(13,11).Value = 6000
$ws.Cells.Item(13,12).Value = 7000
$ws.Cells.Item(13,13).Value = 6500
$ws.Cells.Item(13,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(13,16).Value = 108
$ws.Cells.Item(13,17).Value = 60

# Row 14
$ws.Cells.Item(14,4).Value = 44494
$ws.Cells.Item(14,9).Value = "Primera"
$ws.Cells.Item(14,10).Value = 120
$ws.Cells.Item(14,11).Value = 5000
$ws.Cells.Item(14,12).Value = 6000
$ws.Cells.Item(14,13).Value = 5500
$ws.Cells.Item(14,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(14,16).Value = 92
$ws.Cells.Item(14,17).Value = 60

# Row 15
$ws.Cells.Item(15,4).Value = 45044
$ws.Cells.Item(15,9).Value = "Primera"
$ws.Cells.Item(15,10).Value = 190
$ws.Cells.Item(15,11).Value = 4000
$ws.Cells.Item(15,12).Value = 5000
$ws.Cells.Item(15,13).Value = 4526
$ws.Cells.Item(15,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15,16).Value = 75
$ws.Cells.Item(15,17).Value = 60

# Row 16
$ws.Cells.Item(16,4).Value = 44603
$ws.Cells.Item(16,9).Value = "Primera"
$ws.Cells.Item(16,10).Value = 140
$ws.Cells.Item(16,11).Value = 5500
$ws.Cells.Item(16,12).Value = 6000
$ws.Cells.Item(16,13).Value = 5750
$ws.Cells.Item(16,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(16,16).Value = 96
$ws.Cells.Item(16,17).Value = 60

# Row 17
$ws.Cells.Item(17,4).Value = 44648
$ws.Cells.Item(17,9).Value = "Primera"
$ws.Cells.Item(17,10).Value = 120
$ws.Cells.Item(17,11).Value = 6500
$ws.Cells.Item(17,12).Value = 7000
$ws.Cells.Item(17,13).Value = 6750
$ws.Cells.Item(17,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(17,16).Value = 112
$ws.Cells.Item(17,17).Value = 60

# Row 18
$ws.Cells.Item(18,4).Value = 44967
$ws.Cells.Item(18,9).Value = "Segunda"
$ws.Cells.Item(18,10).Value = 50
$ws.Cells.Item(18,11).Value = 4500
$ws.Cells.Item(18,12).Value = 5000
$ws.Cells.Item(18,13).Value = 4850
$ws.Cells.Item(18,14).Value = "$/caja 90 unidades"
$ws.Cells.Item(18,16).Value = 54
$ws.Cells.Item(18,17).Value = 90

# Row 19
$ws.Cells.Item(19,4).Value = 44421
$ws.Cells.Item(19,9).Value = "Primera"
$ws.Cells.Item(19,10).Value = 100
$ws.Cells.Item(19,11).Value = 8000
$ws.Cells.Item(19,12).Value = 9000
$ws.Cells.Item(19,13).Value = 8500
$ws.Cells.Item(19,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(19,16).Value = 142
$ws.Cells.Item(19,17).Value = 60

# Row 20
$ws.Cells.Item(20,4).Value = 44657
$ws.Cells.Item(20,9).Value = "Primera"
$ws.Cells.Item(20,10).Value = 100
$ws.Cells.Item(20,11).Value = 5000
$ws.Cells.Item(20,12).Value = 5500
$ws.Cells.Item(20,13).Value = 5250
$ws.Cells.Item(20,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(20,16).Value = 88
$ws.Cells.Item(20,17).Value = 60

# Row 21
$ws.Cells.Item(21,4).Value = 44676
$ws.Cells.Item(21,9).Value = "Primera"
$ws.Cells.Item(21,10).Value = 120
$ws.Cells.Item(21,11).Value = 4000
$ws.Cells.Item(21,12).Value = 4500
$ws.Cells.Item(21,13).Value = 4250
$ws.Cells.Item(21,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(21,16).Value = 71
$ws.Cells.Item(21,17).Value = 60

# Row 22
$ws.Cells.Item(22,4).Value = 44935
$ws.Cells.Item(22,9).Value = "Primera"
$ws.Cells.Item(22,10).Value = 120
$ws.Cells.Item(22,11).Value = 6000
$ws.Cells.Item(22,12).Value = 7000
$ws.Cells.Item(22,13).Value = 6500
$ws.Cells.Item(22,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(22,16).Value = 108
$ws.Cells.Item(22,17).Value = 60

# Row 23
$ws.Cells.Item(23,4).Value = 45079
$ws.Cells.Item(23,9).Value = "Primera"
$ws.Cells.Item(23,10).Value = 130
$ws.Cells.Item(23,11).Value = 4000
$ws.Cells.Item(23,12).Value = 5000
$ws.Cells.Item(23,13).Value = 4462
$ws.Cells.Item(23,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(23,16).Value = 74
$ws.Cells.Item(23,17).Value = 60

# Row 24
$ws.Cells.Item(24,4).Value = 44242
$ws.Cells.Item(24,9).Value = "Primera"
$ws.Cells.Item(24,10).Value = 160
$ws.Cells.Item(24,11).Value = 5000
$ws.Cells.Item(24,12).Value = 5500
$ws.Cells.Item(24,13).Value = 5250
$ws.Cells.Item(24,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(24,16).Value = 88
$ws.Cells.Item(24,17).Value = 60

# Row 25
$ws.Cells.Item(25,4).Value = 44589
$ws.Cells.Item(25,9).Value = "Primera"
$ws.Cells.Item(25,10).Value = 110
$ws.Cells.Item(25,11).Value = 5000
$ws.Cells.Item(25,12).Value = 6000
$ws.Cells.Item(25,13).Value = 5500
$ws.Cells.Item(25,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(25,16).Value = 92
$ws.Cells.Item(25,17).Value = 60
